$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gridnodes")

# trafo_id -> gridnode_id refactor: the value "T1" in A2 becomes "T0"
$ws.Range("A2").Value = "T0"

# Update the active selection to match the edited workbook (C8)
$ws.Range("C8").Select()
